# overhauled project, including new pdf to txt conversion
#
# Appends two newly-scraped DEFM14A .htm links to the bottom of the
# "url" column, and converts the first existing link (A2) into a real
# clickable hyperlink (matching what Excel does automatically when it
# recognizes a URL typed into a cell - it mints the built-in "Hyperlink"
# cell style and records the link target).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn the existing A2 URL into a live hyperlink -----------------------
# (applies Excel's built-in "Hyperlink" style: underline + theme color 10)
$existingUrl = $ws.Range("A2").Text
$ws.Hyperlinks.Add($ws.Range("A2"), $existingUrl)

# --- Append the two newly scraped links to the bottom of column A ---------
$ws.Range("A4").Value = "https://www.sec.gov/Archives/edgar/data/1037760/000119312516733208/d263126ddefm14a.htm"
$ws.Range("A5").Value = "https://www.sec.gov/Archives/edgar/data/1226308/000119312516744017/d215313ddefm14a.htm"

# --- Leave the selection where the author last left it ---------------------
$null = $ws.Range("G9").Select()
